$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Estimated Effort" tags in column C for each backlog item.
# Values are written top-to-bottom so new shared-strings get appended in the
# same order the author's Excel session produced them.
$ws.Range("C5").Value  = "Little"
$ws.Range("C6").Value  = "Moderate"
$ws.Range("C7").Value  = "Moderate"
$ws.Range("C8").Value  = "Moderate"
$ws.Range("C9").Value  = "Moderate"
$ws.Range("C10").Value = "Very Heavy"
$ws.Range("C11").Value = "Little"
$ws.Range("C12").Value = "Moderate"

$ws.Range("C15").Value = "Little"

$ws.Range("C17").Value = "Moderate"
$ws.Range("C18").Value = "Little"
$ws.Range("C19").Value = "Little"
$ws.Range("C21").Value = "Moderate"
$ws.Range("C22").Value = "Little"
$ws.Range("C23").Value = "Heavy"
$ws.Range("C24").Value = "Moderate"
$ws.Range("C25").Value = "Moderate"
$ws.Range("C27").Value = "Moderate"
$ws.Range("C28").Value = "Little"
$ws.Range("C29").Value = "Heavy"
$ws.Range("C30").Value = "Moderate"
$ws.Range("C31").Value = "Moderate"

$ws.Range("C34").Value = "Little"
$ws.Range("C35").Value = "Little"

# Move the active cell/selection to where the author left off editing.
$ws.Range("C36").Select()
